$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.759.29'
$ws.Range("E2").Value = '  -3.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.612.13'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.92'
$ws.Range("E5").Value = '  -4.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.15'
$ws.Range("E6").Value = '  -2.45%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.624'
$ws.Range("E8").Value = '  -2.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.608.53'
$ws.Range("E9").Value = '  -1.93%  '
$ws.Range("E10").Value = '  -6.97%  '
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.381'
$ws.Range("E12").Value = '  -4.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.14'
$ws.Range("E14").Value = '  -3.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.083.05'
$ws.Range("E15").Value = '  -1.92%  '
$ws.Range("E16").Value = '  -7.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.640.85'
$ws.Range("E17").Value = '  -3.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.631.53'
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.01'
$ws.Range("E19").Value = '  -4.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.54'
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.52'
$ws.Range("E21").Value = '  -6.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '343.15'
$ws.Range("E22").Value = '  -3.58%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.38'
$ws.Range("E24").Value = '  -3.82%  '
$ws.Range("E25").Value = '  +1.37%  '
$ws.Range("E26").Value = '  -4.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '598.27'
$ws.Range("E27").Value = '  +3.68%  '
$ws.Range("E28").Value = '  -6.08%  '
$ws.Range("E29").Value = '  -3.43%  '
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.93'
$ws.Range("E32").Value = '  -2.93%  '
$ws.Range("E33").Value = '  -3.34%  '
$ws.Range("E34").Value = '  -5.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.61'
$ws.Range("E35").Value = '  -2.06%  '
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("E37").Value = '  -5.18%  '
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.73'
$ws.Range("E39").Value = '  -4.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '154.51'
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  -5.09%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.51'
$ws.Range("E43").Value = '  +1.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.50'
$ws.Range("E44").Value = '  -3.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '157.04'
$ws.Range("E45").Value = '  -3.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.80'
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.89'
$ws.Range("E47").Value = '  -5.51%  '
$ws.Range("E48").Value = '  -4.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.629'
$ws.Range("E49").Value = '  -2.43%  '
$ws.Range("E50").Value = '  -1.48%  '
$ws.Range("E51").Value = '  -4.90%  '
